$d = $word.ActiveDocument

# Locate the "Data Engineering and Infrastructure Architecture" paragraph
# under the Siege Analytics / PARTNER entry, then insert three new bullet
# paragraphs immediately after it (before the existing "Architect
# enterprise-scale..." bullet).
$headingRange = $d.Content
$headingRange.Find.Execute("Data Engineering and Infrastructure Architecture", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$headingPara = $headingRange.Paragraphs.Item(1)

$newLines = @(
    "• Architected data infrastructure processing 15+ billion voter records to support meta-analytical voter file corrections",
    "• Built scalable ETL pipelines enabling analysis of 50,000+ electoral boundaries across all levels of government",
    "• Developed Python boundary estimation algorithm that reduced mapping costs by 75% for 200+ organizations"
)

$insertAfter = $headingPara.Range
foreach ($line in $newLines) {
    $insertAfter.InsertParagraphAfter()
    $newPara = $insertAfter.Next(4, 1).Paragraphs.Item(1)
    $newPara.Range.Text = $line
    $insertAfter = $newPara.Range
}
